$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on every touched cell first so Excel COM
# does not silently re-interpret numeric-looking strings (prices like
# "0.160" / "1.00", percentages, etc.) as native numbers, which would
# drop significant trailing zeros or reformat the text.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.718.61'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -3.06%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.611.18'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.98%  '

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.20%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '574.54'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -3.72%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '155.72'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.85%  '

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.17%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.622'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -5.02%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.119'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -5.52%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.82'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.00%  '

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -4.35%  '

$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.19%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '28.15'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -2.13%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.088.35'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.58%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000182'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -6.55%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.566.21'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -3.06%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.592.29'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.55%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.04'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -4.03%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.64'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +2.92%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.55'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -4.75%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '343.03'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -1.87%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.00'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.07%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.36'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -3.11%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.76'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -1.77%  '

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -2.83%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '589.17'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +4.65%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.19'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -3.50%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.57'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -2.45%  '

$ws.Range('B29').NumberFormat = '@'
$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').NumberFormat = '@'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.31%  '

$ws.Range('B30').NumberFormat = '@'
$ws.Range('B30').Value = 'Kaspa'
$ws.Range('C30').NumberFormat = '@'
$ws.Range('C30').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.160'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -1.19%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.91'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.68%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.06'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -4.69%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.75'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -3.39%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.55'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.25%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.35'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -2.10%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.407'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -3.37%  '

$ws.Range('B37').NumberFormat = '@'
$ws.Range('B37').Value = 'FirstDigitalUSD'
$ws.Range('C37').NumberFormat = '@'
$ws.Range('C37').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.13%  '

$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'EthereumClassic'
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.70'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -3.65%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '154.49'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.20%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.86'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -3.75%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '41.48'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -2.78%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.41'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +6.01%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '155.65'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -3.10%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.90'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -4.01%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '23.28'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +3.00%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0591'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -2.34%  '

$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.628'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.72%  '

$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.102'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.88%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0246'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -3.41%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '18.90'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -4.25%  '
